$d = $word.ActiveDocument
$d.Content.Find.Execute("(d.mine_number)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "({d.mine_number})", 2)
